$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab08")

# Fix mojibake (UTF-8 mis-decoded as Latin-1/Windows-1252) text in the
# "Regional Economic Communities" footnote (cell A103): PaÃ­ses -> Países,
# LÃ­ngua -> Língua, ComÃºn -> Común.
$ws.Range("A103").Value = 'Regional Economic Communities: CEN-SAD = "Community of Sahel-Saharan States"; COMESA = "Common Market for Eastern and Southern Africa"; EAC = "East African Community"; ECCAS = "Economic Community of Central African States"; ECOWAS = "Economic Community of West African States"; IGAD = "Intergovernmental Authority on Development"; SADC = "Southern African Development Community"; UMA = "Arab Maghreb Union"; PALOP = "Países Africanos de Língua Oficial Portuguesa"; ASEAN = "Association of Southeast Asian Nations"; MERCOSUR = "Mercado Común del Sur". EU27 = "European Union (27 members)". OECD = "Organisation for Economic Co-operation and Development".'

# Statistical data refresh: updated aggregate figures for several
# country-grouping rows (regional/income-group rows and sub-region rows).

$ws.Range("C13").Value = 1.5503800000000001
$ws.Range("D13").Value = 1.321766
$ws.Range("E13").Value = 0.337065
$ws.Range("F13").Value = 6.0719950000000003
$ws.Range("G13").Value = 3.833402
$ws.Range("H13").Value = 13.114608
$ws.Range("D23").Value = 5.279325
$ws.Range("E23").Value = 0.68224399999999996
$ws.Range("F23").Value = 7.8306639999999996
$ws.Range("G23").Value = 9.0390119999999996
$ws.Range("H23").Value = 22.831244999999999
$ws.Range("D38").Value = 3.769666
$ws.Range("F38").Value = 0.45433200000000001
$ws.Range("G38").Value = 1.6328659999999999
$ws.Range("H38").Value = 5.9029699999999998
$ws.Range("E45").Value = 3.1045159999999998
$ws.Range("H45").Value = 10.049448
$ws.Range("C61").Value = 0.028899000000000001
$ws.Range("D61").Value = 1.801922
$ws.Range("E61").Value = 0.85586499999999999
$ws.Range("F61").Value = 4.6677530000000003
$ws.Range("G61").Value = 1.5045409999999999
$ws.Range("H61").Value = 8.8508870000000002
$ws.Range("C62").Value = 0.279443
$ws.Range("D62").Value = 1.6026260000000001
$ws.Range("E62").Value = 1.4179409999999999
$ws.Range("F62").Value = 5.1557810000000002
$ws.Range("G62").Value = 1.8272079999999999
$ws.Range("H62").Value = 10.279185999999999
$ws.Range("C63").Value = 0.33656700000000001
$ws.Range("D63").Value = 0.141568
$ws.Range("E63").Value = 0.69418500000000005
$ws.Range("F63").Value = 1.7462930000000001
$ws.Range("G63").Value = 0.93595899999999999
$ws.Range("H63").Value = 3.85433
$ws.Range("C64").Value = 0.067177000000000001
$ws.Range("D64").Value = 0.42281400000000002
$ws.Range("E64").Value = 0.151168
$ws.Range("F64").Value = 2.0261740000000001
$ws.Range("G64").Value = 3.5334159999999999
$ws.Range("H64").Value = 6.2007490000000001
$ws.Range("C65").Value = 0.68862400000000001
$ws.Range("D65").Value = 0.18616199999999999
$ws.Range("E65").Value = 0.62605599999999995
$ws.Range("F65").Value = 1.2552570000000001
$ws.Range("G65").Value = 0.88283299999999998
$ws.Range("H65").Value = 3.6381869999999998
$ws.Range("C66").Value = 0.33339000000000002
$ws.Range("D66").Value = 0.22279199999999999
$ws.Range("E66").Value = 0.73431900000000006
$ws.Range("F66").Value = 1.9358420000000001
$ws.Range("G66").Value = 0.98550599999999999
$ws.Range("H66").Value = 4.2115030000000004
$ws.Range("C67").Value = 0.013698999999999999
$ws.Range("D67").Value = 1.98525
$ws.Range("E67").Value = 1.1929529999999999
$ws.Range("F67").Value = 3.302
$ws.Range("G67").Value = 2.2189459999999999
$ws.Range("H67").Value = 8.7128479999999993
$ws.Range("C68").Value = 0.0041520000000000003
$ws.Range("D68").Value = 1.1138239999999999
$ws.Range("E68").Value = 1.2244120000000001
$ws.Range("F68").Value = 4.6314260000000003
$ws.Range("G68").Value = 1.0230980000000001
$ws.Range("H68").Value = 7.9912640000000001
$ws.Range("C69").Value = 0.012645
$ws.Range("D69").Value = 4.2016470000000004
$ws.Range("E69").Value = 0.059493999999999998
$ws.Range("F69").Value = 0.110413
$ws.Range("G69").Value = 5.9285560000000004
$ws.Range("H69").Value = 10.312754999999999
$ws.Range("D70").Value = 3.40455
$ws.Range("E70").Value = 0.78273000000000004
$ws.Range("F70").Value = 15.522891
$ws.Range("G70").Value = 5.0505639999999996
$ws.Range("H70").Value = 24.760735
$ws.Range("C71").Value = 0.028899000000000001
$ws.Range("D71").Value = 1.801922
$ws.Range("E71").Value = 0.85586499999999999
$ws.Range("F71").Value = 4.6677530000000003
$ws.Range("G71").Value = 1.5045409999999999
$ws.Range("H71").Value = 8.8508870000000002
$ws.Range("D72").Value = 4.115189
$ws.Range("F72").Value = 0.60771799999999998
$ws.Range("G72").Value = 1.312568
$ws.Range("H72").Value = 6.0355309999999998
$ws.Range("C73").Value = 1.2119930000000001
$ws.Range("D73").Value = 2.0912890000000002
$ws.Range("E73").Value = 0.28636099999999998
$ws.Range("F73").Value = 4.7743630000000001
$ws.Range("G73").Value = 5.4736370000000001
$ws.Range("H73").Value = 13.837643
$ws.Range("E74").Value = 4.5813119999999996
$ws.Range("F74").Value = 11.861803
$ws.Range("H74").Value = 16.898354999999999
$ws.Range("C76").Value = 0.53217300000000001
$ws.Range("D76").Value = 0.626641
$ws.Range("E76").Value = 1.0488649999999999
$ws.Range("F76").Value = 0.69576300000000002
$ws.Range("G76").Value = 0.97125499999999998
$ws.Range("H76").Value = 3.871305
$ws.Range("C77").Value = 0.092829999999999996
$ws.Range("D77").Value = 0.52990099999999996
$ws.Range("E77").Value = 0.16214999999999999
$ws.Range("F77").Value = 2.2284039999999998
$ws.Range("G77").Value = 4.5868799999999998
$ws.Range("H77").Value = 7.6001649999999996
$ws.Range("C78").Value = 0.023935000000000001
$ws.Range("D78").Value = 0.062524999999999997
$ws.Range("F78").Value = 0.046670000000000003
$ws.Range("G78").Value = 0.080012
$ws.Range("H78").Value = 0.26608700000000002
$ws.Range("C79").Value = 0.094674999999999995
$ws.Range("D79").Value = 0.053699999999999998
$ws.Range("E79").Value = 0.24251
$ws.Range("F79").Value = 0.50128600000000001
$ws.Range("G79").Value = 0.54537199999999997
$ws.Range("H79").Value = 1.437543
$ws.Range("C80").Value = 0.070668999999999996
$ws.Range("D80").Value = 0.65700899999999995
$ws.Range("E80").Value = 4.9961359999999999
$ws.Range("F80").Value = 21.286928
$ws.Range("G80").Value = 0.218638
$ws.Range("H80").Value = 27.229379000000002
$ws.Range("C81").Value = 0.34847099999999998
$ws.Range("D81").Value = 0.15971099999999999
$ws.Range("E81").Value = 5.133038
$ws.Range("F81").Value = 15.052368
$ws.Range("G81").Value = 2.20017
$ws.Range("H81").Value = 22.893757999999998
$ws.Range("C82").Value = 0.31432399999999999
$ws.Range("D82").Value = 1.7606139999999999
$ws.Range("E82").Value = 0.81823500000000005
$ws.Range("F82").Value = 2.4606849999999998
$ws.Range("G82").Value = 2.0959590000000001
$ws.Range("H82").Value = 7.4472480000000001
$ws.Range("C83").Value = 0.33546100000000001
$ws.Range("D83").Value = 0.13988400000000001
$ws.Range("E83").Value = 0.28206700000000001
$ws.Range("F83").Value = 0.51090800000000003
$ws.Range("G83").Value = 0.81859000000000004
$ws.Range("H83").Value = 2.0867119999999999
$ws.Range("C84").Value = 0.19755300000000001
$ws.Range("D84").Value = 5.9318220000000004
$ws.Range("E84").Value = 0.19320599999999999
$ws.Range("F84").Value = 1.583286
$ws.Range("G84").Value = 7.7685320000000004
$ws.Range("H84").Value = 15.669418
$ws.Range("C86").Value = 0.013327
$ws.Range("D86").Value = 1.039587
$ws.Range("E86").Value = 1.087318
$ws.Range("F86").Value = 4.7706090000000003
$ws.Range("G86").Value = 0.95809599999999995
$ws.Range("H86").Value = 7.8689359999999997
$ws.Range("C87").Value = 0.86368100000000003
$ws.Range("D87").Value = 0.260017
$ws.Range("E87").Value = 0.46596900000000002
$ws.Range("F87").Value = 0.31593300000000002
$ws.Range("G87").Value = 1.2412559999999999
$ws.Range("H87").Value = 3.143942
$ws.Range("C88").Value = 1.2392069999999999
$ws.Range("D88").Value = 0.46784300000000001
$ws.Range("E88").Value = 3.3722379999999998
$ws.Range("F88").Value = 9.2073219999999996
$ws.Range("G88").Value = 1.871783
$ws.Range("H88").Value = 16.158391999999999
$ws.Range("C89").Value = 0.46223900000000001
$ws.Range("D89").Value = 0.196071
$ws.Range("E89").Value = 0.61134200000000005
$ws.Range("F89").Value = 1.917108
$ws.Range("G89").Value = 1.3170040000000001
$ws.Range("H89").Value = 4.503749
$ws.Range("C90").Value = 0.124586
$ws.Range("D90").Value = 0.075142
$ws.Range("E90").Value = 0.80691199999999996
$ws.Range("F90").Value = 1.9680880000000001
$ws.Range("G90").Value = 0.59853400000000001
$ws.Range("H90").Value = 3.5730729999999999
$ws.Range("C91").Value = 0.110486
$ws.Range("D91").Value = 4.3705410000000002
$ws.Range("E91").Value = 0.264733
$ws.Range("F91").Value = 5.0276810000000003
$ws.Range("G91").Value = 5.1283099999999999
$ws.Range("H91").Value = 14.898704
$ws.Range("C92").Value = 0.016351000000000001
$ws.Range("D92").Value = 0.477238
$ws.Range("E92").Value = 0.99463699999999999
$ws.Range("F92").Value = 0.056831
$ws.Range("G92").Value = 0.363871
$ws.Range("H92").Value = 1.9084159999999999
$ws.Range("D94").Value = 0.16491700000000001
$ws.Range("E94").Value = 0.66266599999999998
$ws.Range("F94").Value = 0.61307100000000003
$ws.Range("G94").Value = 0.98255899999999996
$ws.Range("H94").Value = 2.422479
$ws.Range("C95").Value = 0.073090000000000002
$ws.Range("D95").Value = 4.9610279999999998
$ws.Range("E95").Value = 0.0038279999999999998
$ws.Range("F95").Value = 0.80238900000000002
$ws.Range("G95").Value = 4.4098350000000002
$ws.Range("H95").Value = 10.250170000000001
$ws.Range("C97").Value = 0.056774999999999999
$ws.Range("D97").Value = 2.656631
$ws.Range("E97").Value = 0.72939600000000004
$ws.Range("F97").Value = 6.3631120000000001
$ws.Range("G97").Value = 2.466046
$ws.Range("H97").Value = 12.268167999999999
$ws.Range("C98").Value = 0.040826000000000001
$ws.Range("D98").Value = 0.245892
$ws.Range("E98").Value = 2.9416609999999999
$ws.Range("F98").Value = 9.2350999999999992
$ws.Range("G98").Value = 1.2278309999999999
$ws.Range("H98").Value = 13.690860000000001
